$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bulk account number to the masked template format 1234***5678
$ws.Range("B2").Value = "1660***5758"

# Preserve the "text number" quote-prefix formatting used by the sibling
# account-number cells (B3/B4) so the masked value keeps the same style.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats

# Reflect the active selection as it was left after the edit
$ws.Range("B2").Select()
